$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 31   Number  12"
$ws.Range("C9").Value = "Report Covering the Week  3/18/2024  Through  3/24/2024"

# --- Style/type conversions first, using untouched donor cells (rows 33/39 are outside the edited range) ---
# donors: C33 (text "0", style 14), E33 (text "***.*", style 14), I39 (numeric style 15), K39 (numeric style 16)
$ws.Range("C33").Copy($ws.Range("G15"))   # G15: numeric -> text "0"
$ws.Range("E33").Copy($ws.Range("H15"))   # H15: numeric -> text "***.*"
$ws.Range("I39").Copy($ws.Range("C16"))   # C16: text "0" -> numeric
$ws.Range("C33").Copy($ws.Range("C17"))   # C17: numeric -> text "0"
$ws.Range("K39").Copy($ws.Range("L22"))   # L22: text "***.*" -> numeric
$ws.Range("C33").Copy($ws.Range("G27"))   # G27: numeric -> text "0"
$ws.Range("E33").Copy($ws.Range("H27"))   # H27: numeric -> text "***.*"
$ws.Range("I39").Copy($ws.Range("C28"))   # C28: text "0" -> numeric
$ws.Range("C33").Copy($ws.Range("C31"))   # C31: numeric -> text "0"

# --- Final cell values ---
$ws.Range("N14").Value = -50
$ws.Range("L15").Value = 40
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -40
$ws.Range("I16").Value = 19
$ws.Range("J16").Value = 18
$ws.Range("K16").Value = 5.555555555555
$ws.Range("L16").Value = -26.923076923076
$ws.Range("M16").Value = -54.761904761904
$ws.Range("N16").Value = -89.617486338797
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = -17.647058823529
$ws.Range("I17").Value = 43
$ws.Range("J17").Value = 47
$ws.Range("K17").Value = -8.510638297872
$ws.Range("L17").Value = -4.444444444444
$ws.Range("M17").Value = 38.709677419354
$ws.Range("N17").Value = -29.508196721311
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 7
$ws.Range("H18").Value = -22.222222222222
$ws.Range("I18").Value = 25
$ws.Range("J18").Value = 30
$ws.Range("K18").Value = -16.666666666666
$ws.Range("L18").Value = -43.181818181818
$ws.Range("M18").Value = -71.91011235955
$ws.Range("N18").Value = -93.872549019607
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 35
$ws.Range("H19").Value = 25.714285714285
$ws.Range("I19").Value = 129
$ws.Range("J19").Value = 129
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 3.2
$ws.Range("M19").Value = 53.571428571428
$ws.Range("N19").Value = -14
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 66.666666666666
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = -8.333333333333
$ws.Range("I20").Value = 42
$ws.Range("J20").Value = 33
$ws.Range("K20").Value = 27.272727272727
$ws.Range("L20").Value = 121.052631578947
$ws.Range("M20").Value = 23.529411764705
$ws.Range("N20").Value = -90.410958904109
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = -28
$ws.Range("F21").Value = 83
$ws.Range("G21").Value = 83
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 267
$ws.Range("J21").Value = 261
$ws.Range("K21").Value = 2.298850574712
$ws.Range("L21").Value = 1.136363636363
$ws.Range("M21").Value = -4.982206405693
$ws.Range("N21").Value = -78.588612670409
$ws.Range("L22").Value = 100
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = -50
$ws.Range("F24").Value = 69
$ws.Range("G24").Value = 80
$ws.Range("H24").Value = -13.75
$ws.Range("I24").Value = 217
$ws.Range("J24").Value = 252
$ws.Range("K24").Value = -13.888888888888
$ws.Range("L24").Value = -11.065573770491
$ws.Range("M24").Value = 10.714285714285
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -50
$ws.Range("I25").Value = 53
$ws.Range("J25").Value = 74
$ws.Range("K25").Value = -28.378378378378
$ws.Range("L25").Value = -24.285714285714
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 38
$ws.Range("H26").Value = 58.333333333333
$ws.Range("I26").Value = 99
$ws.Range("J26").Value = 63
$ws.Range("K26").Value = 57.142857142857
$ws.Range("L26").Value = 39.436619718309
$ws.Range("M26").Value = 8.791208791208
$ws.Range("L27").Value = 42.857142857142
$ws.Range("C28").Value = 2
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 8
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 166.666666666667
$ws.Range("I28").Value = 22
$ws.Range("J28").Value = 8
$ws.Range("K28").Value = 175
$ws.Range("L28").Value = 69.230769230769
$ws.Range("L31").Value = 0

Write-Output "done"